$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 2..11: column A is ball number (numeric), column B is distance (stored as text)
$ballNumbers = @(1, 2, 3, 4, 5, 6, 7, 8, 9, 10)
$distances = @("0.040041", "0.061339", "0.04219", "0.137089", "0.101683", "0.056507", "0.05839", "0.040633", "0.037811", "0.016227")

for ($i = 0; $i -lt $ballNumbers.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $ballNumbers[$i]
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = $distances[$i]
}
